# Update res_bus/vm_pu.xlsx values for the 380 kV case (Case_2_194).
# Rows 2-25 (bus index 0-23), columns B-F and I-N get new per-unit voltage values;
# the base-case vm_pu at the slack bus (col B) moves from 1.05 to 1.02 pu.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046919633519343
$ws.Range("D2").Value = 1.056451915327412
$ws.Range("E2").Value = 1.060710624707767
$ws.Range("F2").Value = 1.068653181176422
$ws.Range("I2").Value = 1.049859609129481
$ws.Range("J2").Value = 1.051971462113002
$ws.Range("K2").Value = 1.059189588494537
$ws.Range("L2").Value = 1.063436662836598
$ws.Range("M2").Value = 1.071357787801804
$ws.Range("N2").Value = 1.053465382252016

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047790816784484
$ws.Range("D3").Value = 1.057165926946021
$ws.Range("E3").Value = 1.061600813002955
$ws.Range("F3").Value = 1.069517547882628
$ws.Range("I3").Value = 1.050121720316649
$ws.Range("J3").Value = 1.052491294062492
$ws.Range("K3").Value = 1.059717495514553
$ws.Range("L3").Value = 1.064141138791813
$ws.Range("M3").Value = 1.072038053710377
$ws.Range("N3").Value = 1.053985952422499

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048354973520782
$ws.Range("D4").Value = 1.057628327024214
$ws.Range("E4").Value = 1.062178061612065
$ws.Range("F4").Value = 1.070077747200889
$ws.Range("I4").Value = 1.050290273244361
$ws.Range("J4").Value = 1.052827433100125
$ws.Range("K4").Value = 1.060058798869752
$ws.Range("L4").Value = 1.064597578745276
$ws.Range("M4").Value = 1.072478481261323
$ws.Range("N4").Value = 1.054322568816122

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048592249654525
$ws.Range("D5").Value = 1.057822811015639
$ws.Range("E5").Value = 1.062421031408992
$ws.Range("F5").Value = 1.070313467837267
$ws.Range("I5").Value = 1.050360880648621
$ws.Range("J5").Value = 1.052968690480304
$ws.Range("K5").Value = 1.060202212690216
$ws.Range("L5").Value = 1.064789607562575
$ws.Range("M5").Value = 1.072663695420058
$ws.Range("N5").Value = 1.054464026797983

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048632095457307
$ws.Range("D6").Value = 1.057855471023063
$ws.Range("E6").Value = 1.062461844325575
$ws.Range("F6").Value = 1.070353058800455
$ws.Range("I6").Value = 1.050372721128981
$ws.Range("J6").Value = 1.052992404930351
$ws.Range("K6").Value = 1.060226288355178
$ws.Range("L6").Value = 1.064821858299486
$ws.Range("M6").Value = 1.072694797076038
$ws.Range("N6").Value = 1.054487774925269

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04835814360666
$ws.Range("D7").Value = 1.05763092537328
$ws.Range("E7").Value = 1.062181307031428
$ws.Range("F7").Value = 1.070080896077403
$ws.Range("I7").Value = 1.050291217696038
$ws.Range("J7").Value = 1.052829320807529
$ws.Range("K7").Value = 1.060060715449244
$ws.Range("L7").Value = 1.06460014408969
$ws.Range("M7").Value = 1.072480955874244
$ws.Range("N7").Value = 1.054324459204287

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047213961437766
$ws.Range("D8").Value = 1.056693137838632
$ws.Range("E8").Value = 1.061011210793207
$ws.Range("F8").Value = 1.068945111338482
$ws.Range("I8").Value = 1.049948407904159
$ws.Range("J8").Value = 1.052147188076068
$ws.Range("K8").Value = 1.059368056002134
$ws.Range("L8").Value = 1.063674619325699
$ws.Range("M8").Value = 1.071587634216876
$ws.Range("N8").Value = 1.053641357766105

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045201219971748
$ws.Range("D9").Value = 1.055043669887546
$ws.Range("E9").Value = 1.058958903719651
$ws.Range("F9").Value = 1.066950654010352
$ws.Range("I9").Value = 1.049336320388941
$ws.Range("J9").Value = 1.050943494419603
$ws.Range("K9").Value = 1.058145351076226
$ws.Range("L9").Value = 1.062048361653748
$ws.Range("M9").Value = 1.070015465797801
$ws.Range("N9").Value = 1.052435954726583

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043861796727591
$ws.Range("D10").Value = 1.053946159560758
$ws.Range("E10").Value = 1.05759722364038
$ws.Range("F10").Value = 1.065625784685666
$ws.Range("I10").Value = 1.048922918910073
$ws.Range("J10").Value = 1.050139963557329
$ws.Range("K10").Value = 1.05732884060372
$ws.Range("L10").Value = 1.060967391952088
$ws.Range("M10").Value = 1.068968769255792
$ws.Range("N10").Value = 1.051631282758318

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043282399112264
$ws.Range("D11").Value = 1.053471451909819
$ws.Range("E11").Value = 1.057009170242346
$ws.Range("F11").Value = 1.065053254299996
$ws.Range("I11").Value = 1.048742653428249
$ws.Range("J11").Value = 1.049791786408097
$ws.Range("K11").Value = 1.056974971781105
$ws.Range("L11").Value = 1.060500097628198
$ws.Range("M11").Value = 1.06851589184283
$ws.Range("N11").Value = 1.051282611157597

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04306727375545
$ws.Range("D12").Value = 1.053295204318932
$ws.Range("E12").Value = 1.056790977598793
$ws.Range("F12").Value = 1.064840765044198
$ws.Range("I12").Value = 1.048675506233081
$ws.Range("J12").Value = 1.049662422515696
$ws.Range("K12").Value = 1.056843483106458
$ws.Range("L12").Value = 1.060326641243426
$ws.Range("M12").Value = 1.068347726889769
$ws.Range("N12").Value = 1.051153063553631

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043113414838433
$ws.Range("D13").Value = 1.053333006380976
$ws.Range("E13").Value = 1.056837769928213
$ws.Range("F13").Value = 1.064886336808239
$ws.Range("I13").Value = 1.048689918055957
$ws.Range("J13").Value = 1.049690173093706
$ws.Range("K13").Value = 1.056871689937602
$ws.Range("L13").Value = 1.060363842884919
$ws.Range("M13").Value = 1.06838379638171
$ws.Range("N13").Value = 1.051180853540645

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043264614965558
$ws.Range("D14").Value = 1.053456881588376
$ws.Range("E14").Value = 1.056991129534316
$ws.Range("F14").Value = 1.065035686311827
$ws.Range("I14").Value = 1.048737106862732
$ws.Range("J14").Value = 1.049781093864058
$ws.Range("K14").Value = 1.056964103818474
$ws.Range("L14").Value = 1.060485757263314
$ws.Range("M14").Value = 1.068501990167864
$ws.Range("N14").Value = 1.05127190342892

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043357786150282
$ws.Range("D15").Value = 1.053533215836547
$ws.Range("E15").Value = 1.057085650849205
$ws.Range("F15").Value = 1.065127728581464
$ws.Range("I15").Value = 1.048766156478348
$ws.Range("J15").Value = 1.049837108494235
$ws.Range("K15").Value = 1.057021036995876
$ws.Range("L15").Value = 1.060560888353096
$ws.Range("M15").Value = 1.068574820443284
$ws.Range("N15").Value = 1.051327997606296

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043900261805827
$ws.Range("D16").Value = 1.053977675470153
$ws.Range("E16").Value = 1.057636283888051
$ws.Range("F16").Value = 1.06566380595078
$ws.Range("I16").Value = 1.048934856052792
$ws.Range("J16").Value = 1.050163065903279
$ws.Range("K16").Value = 1.057352319190156
$ws.Range("L16").Value = 1.060998421140418
$ws.Range("M16").Value = 1.068998832715677
$ws.Range("N16").Value = 1.051654417912249

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044240699087539
$ws.Range("D17").Value = 1.05425661402853
$ws.Range("E17").Value = 1.057982101149318
$ws.Range("F17").Value = 1.066000381367799
$ws.Range("I17").Value = 1.049040340031012
$ws.Range("J17").Value = 1.050367466103228
$ws.Range("K17").Value = 1.057560040429301
$ws.Range("L17").Value = 1.06127308178962
$ws.Range("M17").Value = 1.069264898887462
$ws.Range("N17").Value = 1.051859108383927

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044439326408461
$ws.Range("D18").Value = 1.054419364385836
$ws.Range("E18").Value = 1.058183961203185
$ws.Range("F18").Value = 1.066196810645124
$ws.Range("I18").Value = 1.049101745439152
$ws.Range("J18").Value = 1.050486665762883
$ws.Range("K18").Value = 1.05768117029837
$ws.Range("L18").Value = 1.06143336115582
$ws.Range("M18").Value = 1.0694201243705
$ws.Range("N18").Value = 1.051978477320771

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044507062651779
$ws.Range("D19").Value = 1.054474866520312
$ws.Range("E19").Value = 1.058252815772804
$ws.Range("F19").Value = 1.066263806647811
$ws.Range("I19").Value = 1.049122662447178
$ws.Range("J19").Value = 1.050527305743513
$ws.Range("K19").Value = 1.057722467258068
$ws.Range("L19").Value = 1.061488024862845
$ws.Range("M19").Value = 1.069473057902056
$ws.Range("N19").Value = 1.052019175014835

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044204167585479
$ws.Range("D20").Value = 1.054226681355428
$ws.Range("E20").Value = 1.057944982630203
$ws.Range("F20").Value = 1.06596425856563
$ws.Range("I20").Value = 1.049029035172117
$ws.Range("J20").Value = 1.050345538320845
$ws.Range("K20").Value = 1.05753775704355
$ws.Range("L20").Value = 1.061243605589511
$ws.Range("M20").Value = 1.069236349020883
$ws.Range("N20").Value = 1.051837149461578

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043220087836825
$ws.Range("D21").Value = 1.053420401206913
$ws.Range("E21").Value = 1.056945962408202
$ws.Range("F21").Value = 1.06499170178771
$ws.Range("I21").Value = 1.048723216134927
$ws.Range("J21").Value = 1.049754320931538
$ws.Range("K21").Value = 1.056936891495166
$ws.Range("L21").Value = 1.060449853261402
$ws.Range("M21").Value = 1.068467183544798
$ws.Range("N21").Value = 1.051245092475765

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04260187179227
$ws.Range("D22").Value = 1.052913924336444
$ws.Range("E22").Value = 1.056319208335283
$ws.Range("F22").Value = 1.064381224432251
$ws.Range("I22").Value = 1.048529844924366
$ws.Range("J22").Value = 1.049382394849576
$ws.Range("K22").Value = 1.056558837708699
$ws.Range("L22").Value = 1.059951470397073
$ws.Range("M22").Value = 1.067983891140277
$ws.Range("N22").Value = 1.050872638216106

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042929550645015
$ws.Range("D23").Value = 1.053182372840909
$ws.Range("E23").Value = 1.056651332125472
$ws.Range("F23").Value = 1.064704753933151
$ws.Range("I23").Value = 1.048632457801594
$ws.Range("J23").Value = 1.049579578892077
$ws.Range("K23").Value = 1.056759275965666
$ws.Range("L23").Value = 1.060215607575065
$ws.Range("M23").Value = 1.068240063491332
$ws.Range("N23").Value = 1.051070102282564

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044220674425842
$ws.Range("D24").Value = 1.054240206486998
$ws.Range("E24").Value = 1.057961754425826
$ws.Range("F24").Value = 1.065980580563724
$ws.Range("I24").Value = 1.049034143727074
$ws.Range("J24").Value = 1.050355446614623
$ws.Range("K24").Value = 1.057547826040504
$ws.Range("L24").Value = 1.061256924385663
$ws.Range("M24").Value = 1.069249249373063
$ws.Range("N24").Value = 1.051847071826269

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04572114413401
$ws.Range("D25").Value = 1.055469727569922
$ws.Range("E25").Value = 1.059488331781911
$ws.Range("F25").Value = 1.067465435993199
$ws.Range("I25").Value = 1.049495504695726
$ws.Range("J25").Value = 1.051254871302107
$ws.Range("K25").Value = 1.058461697033967
$ws.Range("L25").Value = 1.062468230127613
$ws.Range("M25").Value = 1.070421665750591
$ws.Range("N25").Value = 1.052747773799975
